$d = $word.ActiveDocument

# Locate the end of the paragraph that finishes with "most employers." -
# that's where the new paragraph about the collaboration project needs
# to be inserted.
$found = $d.Content
$ok = $found.Find.Execute("most employers.", $true, $false, $false, $false,
                           $false, $true, 1, $false, "", 0)

# Build a clean zero-length insertion point right after the matched text
# (i.e. right before the paragraph mark that ends that paragraph).
$insertPoint = $d.Range($found.End, $found.End)

# Insert a brand-new paragraph (same first-line indent style as its
# neighbours) made up of two runs, matching the target markup exactly.
$run1 = "This collaboration project with the team has also been a great time to sharpen my communicational skills and show what kind of team member that I’d like to become. Within our meetings, I"
$run2 = " keep my topics as short as possible and to the point. While outside the meetings, the communication shifts to a lighthearted and fun demeanor and asking for help or input when needed."

$xmlFragment = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
               "<w:pPr><w:ind w:firstLine='720'/></w:pPr>" +
               "<w:r><w:t>" + $run1 + "</w:t></w:r>" +
               "<w:r><w:t xml:space='preserve'>" + $run2 + "</w:t></w:r>" +
               "</w:p>"

$insertPoint.InsertXML($xmlFragment) | Out-Null
